$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source sheet stores every data value as a shared string (t="s"),
# even the numeric-looking ones. Excel's default behavior for a bare
# ".Value = '39.1'" assignment is to coerce it to a real number, so we
# temporarily force Text number format on the target range, write the
# values as strings, and then clear the formatting again so the cells
# end up with no explicit style (matching the original workbook) while
# still being stored as shared strings.
$dataRange = $ws.Range("A2:H3")
$dataRange.NumberFormat = "@"

# Row 2: update the existing prediction row with the new values
$ws.Range("A2").Value = "39.1"
$ws.Range("B2").Value = "181.0"
$ws.Range("C2").Value = "Adelie"
$ws.Range("D2").Value = "1.0"
$ws.Range("E2").Value = "0.0"
$ws.Range("F2").Value = "0.0"
$ws.Range("G2").Value = "v1.0"
$ws.Range("H2").Value = "2025-05-04 21:10:32"

# Row 3: new second prediction row
$ws.Range("A3").Value = "46.5"
$ws.Range("B3").Value = "192.0"
$ws.Range("C3").Value = "Chinstrap"
$ws.Range("D3").Value = "0.09"
$ws.Range("E3").Value = "0.91"
$ws.Range("F3").Value = "0.0"
$ws.Range("G3").Value = "v1.0"
$ws.Range("H3").Value = "2025-05-04 21:10:32"

# Drop the explicit Text formatting again so the new cells stay unstyled,
# same as the rest of the data rows in the original file.
$dataRange.ClearFormats()
